$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.199.29"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.831.48"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.10"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6187"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07350"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2905"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.21"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.837.22"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.970"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.48"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008936"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.843"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "29.186.87"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "2.083.52"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.17"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.348"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.11"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1394"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.542"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.490"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05829"
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.097"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.849"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7270"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.609"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.857"
$ws.Range("E38").Value = "  +3.14%  "
$ws.Range("D39").Value = "1.221.11"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01752"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.232"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9046"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.989.00"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.72"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.36"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5034"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  -4.45%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4024"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.104"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1131"
$ws.Range("E51").Value = "  +2.83%  "
